$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row by scanning column H (the existing "IP" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row

# --- Header row (row 1): copy the existing header formatting from H1 onto I1/J1 ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows: I = constant 1, J = copy of column H ---
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
